# Apply timesheet update: add two new rows (20 and 21) of work log entries,
# extend the SUM formula in F3 to cover the new rows, and update the
# selected cell in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the existing date style used by the other rows in column A
# (copy format from A19 so no duplicate style entries are introduced)
$ws.Range("A19").Copy()
$ws.Range("A20:A21").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# New row 20: 2016-12-28, domain "DEV", new task description, 2 hours
$ws.Range("A20").Value2 = 42732
$ws.Range("B20").Value = "DEV"
$ws.Range("C20").Value = "User navigatie enkel op de userpagina's laten zien "
$ws.Range("D20").Value = 2

# New row 21: 2016-12-29, domain "DEV", new task description, 3 hours
$ws.Range("A21").Value2 = 42733
$ws.Range("B21").Value = "DEV"
$ws.Range("C21").Value = "Slideshow op iedere aanbod/vraagpagina"
$ws.Range("D21").Value = 3

# Extend the totals formula to include the two newly added rows
$ws.Range("F3").Formula = "=SUM(D2:D21)"

# Update the active selection to reflect where the author last clicked
$ws.Range("F24").Select()
